# Auto-generated Excel COM-interop script to apply the cryptos.xlsx update
# (symbol list refresh from 5-1-2023 to 6-1-2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns D/E/F/G hold numeric-/date-looking text (e.g. "257.72", "-0.09%",
# "6-1-2023", "0") that must stay plain text, matching the original inline-string
# cells. Pre-formatting the whole block as Text ("@") stops Excel's COM layer from
# auto-coercing these into numbers/dates/percentages; re-applying the "Normal" style
# afterwards drops the temporary format so the cells end up styled exactly as before.
$textRange = $ws.Range("D2:G51")
$textRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '257.72'
$ws.Range("E2").Value = '-0.09%'
$ws.Range("F2").Value = '6-1-2023'
$ws.Range("G2").Value = '0'
# Row 3
$ws.Range("D3").Value = '27.06'
$ws.Range("E3").Value = '-1.82%'
$ws.Range("F3").Value = '6-1-2023'
$ws.Range("G3").Value = '0'
# Row 4
$ws.Range("D4").Value = '4.659'
$ws.Range("E4").Value = '-10.84%'
$ws.Range("F4").Value = '6-1-2023'
$ws.Range("G4").Value = '0'
# Row 5
$ws.Range("D5").Value = '0.05899'
$ws.Range("E5").Value = '-0.62%'
$ws.Range("F5").Value = '6-1-2023'
$ws.Range("G5").Value = '0'
# Row 6
$ws.Range("D6").Value = '6.647'
$ws.Range("E6").Value = '-1.20%'
$ws.Range("F6").Value = '6-1-2023'
$ws.Range("G6").Value = '0'
# Row 7
$ws.Range("D7").Value = '0.8567'
$ws.Range("E7").Value = '-1.70%'
$ws.Range("F7").Value = '6-1-2023'
$ws.Range("G7").Value = '0'
# Row 8
$ws.Range("D8").Value = '0.9467'
$ws.Range("E8").Value = '-6.86%'
$ws.Range("F8").Value = '6-1-2023'
$ws.Range("G8").Value = '0'
# Row 9
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").Value = '0.0006027'
$ws.Range("E9").Value = '-0.57%'
$ws.Range("F9").Value = '6-1-2023'
$ws.Range("G9").Value = '0'
# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1407'
$ws.Range("E10").Value = '-1.03%'
$ws.Range("F10").Value = '6-1-2023'
$ws.Range("G10").Value = '0'
# Row 11
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '0.04339'
$ws.Range("E11").Value = '21.96%'
$ws.Range("F11").Value = '6-1-2023'
$ws.Range("G11").Value = '0'
# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.07101'
$ws.Range("E12").Value = '-1.18%'
$ws.Range("F12").Value = '6-1-2023'
$ws.Range("G12").Value = '0'
# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03134'
$ws.Range("E13").Value = '-0.66%'
$ws.Range("F13").Value = '6-1-2023'
$ws.Range("G13").Value = '0'
# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09147'
$ws.Range("E14").Value = '-0.97%'
$ws.Range("F14").Value = '6-1-2023'
$ws.Range("G14").Value = '0'
# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001538'
$ws.Range("E15").Value = '-0.03%'
$ws.Range("F15").Value = '6-1-2023'
$ws.Range("G15").Value = '0'
# Row 16
$ws.Range("D16").Value = '0.005979'
$ws.Range("E16").Value = '4.00%'
$ws.Range("F16").Value = '6-1-2023'
$ws.Range("G16").Value = '0'
# Row 17
$ws.Range("D17").Value = '3.516'
$ws.Range("E17").Value = '0.14%'
$ws.Range("F17").Value = '6-1-2023'
$ws.Range("G17").Value = '0'
# Row 18
$ws.Range("D18").Value = '3.198'
$ws.Range("E18").Value = '-2.03%'
$ws.Range("F18").Value = '6-1-2023'
$ws.Range("G18").Value = '0'
# Row 19
$ws.Range("D19").Value = '2.226'
$ws.Range("E19").Value = '0.80%'
$ws.Range("F19").Value = '6-1-2023'
$ws.Range("G19").Value = '0'
# Row 20
$ws.Range("D20").Value = '0.3055'
$ws.Range("E20").Value = '-3.00%'
$ws.Range("F20").Value = '6-1-2023'
$ws.Range("G20").Value = '0'
# Row 21
$ws.Range("E21").Value = '-0.56%'
$ws.Range("F21").Value = '6-1-2023'
$ws.Range("G21").Value = '0'
# Row 22
$ws.Range("D22").Value = '3.816'
$ws.Range("E22").Value = '7.23%'
$ws.Range("F22").Value = '6-1-2023'
$ws.Range("G22").Value = '0'
# Row 23
$ws.Range("D23").Value = '0.04254'
$ws.Range("E23").Value = '1.09%'
$ws.Range("F23").Value = '6-1-2023'
$ws.Range("G23").Value = '0'
# Row 24
$ws.Range("D24").Value = '0.001220'
$ws.Range("E24").Value = '0.07%'
$ws.Range("F24").Value = '6-1-2023'
$ws.Range("G24").Value = '0'
# Row 25
$ws.Range("D25").Value = '0.004297'
$ws.Range("E25").Value = '-4.75%'
$ws.Range("F25").Value = '6-1-2023'
$ws.Range("G25").Value = '0'
# Row 26
$ws.Range("D26").Value = '0.0001199'
$ws.Range("E26").Value = '-0.08%'
$ws.Range("F26").Value = '6-1-2023'
$ws.Range("G26").Value = '0'
# Row 27
$ws.Range("D27").Value = '0.0001936'
$ws.Range("E27").Value = '-0.13%'
$ws.Range("F27").Value = '6-1-2023'
$ws.Range("G27").Value = '0'
# Row 28
$ws.Range("F28").Value = '6-1-2023'
$ws.Range("G28").Value = '0'
# Row 29
$ws.Range("F29").Value = '6-1-2023'
$ws.Range("G29").Value = '0'
# Row 30
$ws.Range("F30").Value = '6-1-2023'
$ws.Range("G30").Value = '0'
# Row 31
$ws.Range("F31").Value = '6-1-2023'
$ws.Range("G31").Value = '0'
# Row 32
$ws.Range("F32").Value = '6-1-2023'
$ws.Range("G32").Value = '0'
# Row 33
$ws.Range("F33").Value = '6-1-2023'
$ws.Range("G33").Value = '0'
# Row 34
$ws.Range("F34").Value = '6-1-2023'
$ws.Range("G34").Value = '0'
# Row 35
$ws.Range("F35").Value = '6-1-2023'
$ws.Range("G35").Value = '0'
# Row 36
$ws.Range("F36").Value = '6-1-2023'
$ws.Range("G36").Value = '0'
# Row 37
$ws.Range("F37").Value = '6-1-2023'
$ws.Range("G37").Value = '0'
# Row 38
$ws.Range("F38").Value = '6-1-2023'
$ws.Range("G38").Value = '0'
# Row 39
$ws.Range("F39").Value = '6-1-2023'
$ws.Range("G39").Value = '0'
# Row 40
$ws.Range("D40").Value = '0.03824'
$ws.Range("E40").Value = '-0.36%'
$ws.Range("F40").Value = '6-1-2023'
$ws.Range("G40").Value = '0'
# Row 41
$ws.Range("D41").Value = '0.006237'
$ws.Range("E41").Value = '57.59%'
$ws.Range("F41").Value = '6-1-2023'
$ws.Range("G41").Value = '0'
# Row 42
$ws.Range("D42").Value = '0.1105'
$ws.Range("E42").Value = '0.01%'
$ws.Range("F42").Value = '6-1-2023'
$ws.Range("G42").Value = '0'
# Row 43
$ws.Range("D43").Value = '0.002199'
$ws.Range("E43").Value = '-4.84%'
$ws.Range("F43").Value = '6-1-2023'
$ws.Range("G43").Value = '0'
# Row 44
$ws.Range("D44").Value = '0.01145'
$ws.Range("E44").Value = '14.90%'
$ws.Range("F44").Value = '6-1-2023'
$ws.Range("G44").Value = '0'
# Row 45
$ws.Range("D45").Value = '0.00005458'
$ws.Range("E45").Value = '0.15%'
$ws.Range("F45").Value = '6-1-2023'
$ws.Range("G45").Value = '0'
# Row 46
$ws.Range("E46").Value = '-0.07%'
$ws.Range("F46").Value = '6-1-2023'
$ws.Range("G46").Value = '0'
# Row 47
$ws.Range("D47").Value = '0.06164'
$ws.Range("E47").Value = '-43.51%'
$ws.Range("F47").Value = '6-1-2023'
$ws.Range("G47").Value = '0'
# Row 48
$ws.Range("D48").Value = '0.2244'
$ws.Range("E48").Value = '9,941.58%'
$ws.Range("F48").Value = '6-1-2023'
$ws.Range("G48").Value = '0'
# Row 49
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").Value = '-0.07%'
$ws.Range("F49").Value = '6-1-2023'
$ws.Range("G49").Value = '0'
# Row 50
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").Value = '-0.07%'
$ws.Range("F50").Value = '6-1-2023'
$ws.Range("G50").Value = '0'
# Row 51
$ws.Range("F51").Value = '6-1-2023'
$ws.Range("G51").Value = '0'

# Restore default styling on the text-coerced block (drops the temporary "@" format)
$textRange.Style = "Normal"

